$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'71.908.21"
$ws.Range('E2').Value = '  -1.38%  '

$ws.Range('D3').Value = "'2.682.88"
$ws.Range('E3').Value = '  +0.52%  '

$ws.Range('E4').Value = '  +0.05%  '

$ws.Range('D5').Value = "'598.43"
$ws.Range('E5').Value = '  -2.05%  '

$ws.Range('D6').Value = "'175.14"
$ws.Range('E6').Value = '  -3.73%  '

$ws.Range('E7').Value = '  +0.06%  '

$ws.Range('D8').Value = "'0.523"
$ws.Range('E8').Value = '  -1.33%  '

$ws.Range('D9').Value = "'2.681.16"
$ws.Range('E9').Value = '  +0.47%  '

$ws.Range('E10').Value = '  -6.02%  '

$ws.Range('E11').Value = '  +1.91%  '

$ws.Range('E12').Value = '  +0.47%  '

$ws.Range('E13').Value = '  -2.58%  '

$ws.Range('D14').Value = "'3.170.59"
$ws.Range('E14').Value = '  +2.11%  '

$ws.Range('B15').Value = 'ShibaInu'
$ws.Range('C15').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D15').Value = "'0.0000184"
$ws.Range('E15').Value = '  -5.37%  '

$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = "'71.809.78"
$ws.Range('E16').Value = '  -1.29%  '

$ws.Range('D17').Value = "'26.20"
$ws.Range('E17').Value = '  -3.13%  '

$ws.Range('D18').Value = "'2.679.81"
$ws.Range('E18').Value = '  +1.01%  '

$ws.Range('D19').Value = "'12.20"
$ws.Range('E19').Value = '  +4.69%  '

$ws.Range('D20').Value = "'8.23"
$ws.Range('E20').Value = '  +3.35%  '

$ws.Range('D21').Value = "'371.30"
$ws.Range('E21').Value = '  -3.82%  '

$ws.Range('D22').Value = "'4.18"
$ws.Range('E22').Value = '  -1.43%  '

$ws.Range('D23').Value = "'2.01"
$ws.Range('E23').Value = '  -1.90%  '

$ws.Range('D24').Value = "'72.11"
$ws.Range('E24').Value = '  -2.04%  '

$ws.Range('D25').Value = "'0.999"
$ws.Range('E25').Value = '  -0.20%  '

$ws.Range('D26').Value = "'4.33"
$ws.Range('E26').Value = '  -3.61%  '

$ws.Range('E27').Value = '  -2.76%  '

$ws.Range('D28').Value = "'2.816.71"
$ws.Range('E28').Value = '  +0.48%  '

$ws.Range('D29').Value = "'0.998"
$ws.Range('E29').Value = '  -0.23%  '

$ws.Range('D30').Value = "'0.0₃0970"
$ws.Range('E30').Value = '  -1.28%  '

$ws.Range('D31').Value = "'8.04"
$ws.Range('E31').Value = '  -1.44%  '

$ws.Range('D32').Value = "'502.80"
$ws.Range('E32').Value = '  -8.55%  '

$ws.Range('D33').Value = "'1.30"
$ws.Range('E33').Value = '  -3.76%  '

$ws.Range('E34').Value = '  -2.16%  '

$ws.Range('E35').Value = '  +0.02%  '

$ws.Range('D36').Value = "'163.40"
$ws.Range('E36').Value = '  -1.06%  '

$ws.Range('D37').Value = "'19.56"
$ws.Range('E37').Value = '  +0.45%  '

$ws.Range('D38').Value = "'19.08"
$ws.Range('E38').Value = '  -0.41%  '

$ws.Range('D39').Value = "'1.38"
$ws.Range('E39').Value = '  -3.83%  '

$ws.Range('E40').Value = '  -4.49%  '

$ws.Range('E41').Value = '  -5.16%  '

$ws.Range('E42').Value = '  -0.18%  '

$ws.Range('E43').Value = '  -2.87%  '

$ws.Range('D44').Value = "'2.56"
$ws.Range('E44').Value = '  -3.56%  '

$ws.Range('D45').Value = "'0.333"
$ws.Range('E45').Value = '  -1.45%  '

$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = "'157.01"
$ws.Range('E46').Value = '  +2.54%  '

$ws.Range('B47').Value = 'OKB'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D47').Value = "'39.50"
$ws.Range('E47').Value = '  -0.62%  '

$ws.Range('E48').Value = '  +2.80%  '

$ws.Range('D49').Value = "'3.72"
$ws.Range('E49').Value = '  -0.37%  '

$ws.Range('E50').Value = '  +1.20%  '

$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').Value = "'0.608"
$ws.Range('E51').Value = '  -0.33%  '
